# Update the "Förändrad" (Changed) date column (C) for rows 2-8.
# The stored serial date value increases from 46062 to 46063 (one day later)
# for every data row in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 3).Value = 46063
}
